$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 and add new rows 5-6 with the refreshed partner list.
$ws.Range("A2").Value = "Clinique Anne d'Artois"
$ws.Range("B2").Value = "vivalto_sante.jpg"
$ws.Range("C2").Value = "image/jpg"

$ws.Range("A3").Value = "Clinique Saint-Amé"
$ws.Range("B3").Value = "ramsay_sante.png"
$ws.Range("C3").Value = "image/png"

$ws.Range("A4").Value = "Polyclinique Vauban"
$ws.Range("B4").Value = "elsan.png"
$ws.Range("C4").Value = "image/png"

$ws.Range("A5").Value = "SOS Mains Côte d'Opale"
$ws.Range("B5").Value = "sos_mains_cote_opale.jpg"
$ws.Range("C5").Value = "image/jpg"

$ws.Range("A6").Value = "Hôpital Chantilly Les Jockeys"
$ws.Range("B6").Value = "hopital_chantilly_les_jockeys.png"
$ws.Range("C6").Value = "image/png"

# Grow the table (ListObject) to cover the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C6")) | Out-Null

# Widen column A to fit the longer partner names.
$ws.Columns.Item(1).ColumnWidth = 20.5

# Match the author's final selection/cursor position.
$ws.Range("B8").Select() | Out-Null
